# Applies the edit described by the commit "Modificación de los cálculos con altura antena":
#   - Antenna + Headend row 49 (antenna/support heights) now derive the G:N
#     heights from the Mast + Tower sheet instead of being hard-coded, and the
#     first column (F49) becomes a fixed 20 m value.
#   - Fixes a copy/paste range typo in the S/I-with-mast-amplifier formula
#     (J96) so it sums H86:N86 like its neighbours instead of I86:O86.
#   - Clears a couple of leftover Spanish scratch notes (P5, F48 on
#     "Antenna + Headend") and a stray "a" placeholder (L16 on "Mast + Tower")
#     that were no longer needed.

$wb = $excel.ActiveWorkbook

$headend = $wb.Worksheets.Item("Antenna + Headend")
$mast = $wb.Worksheets.Item("Mast + Tower")

# --- Row 49: antenna support height per outlet -----------------------------
# F49 becomes a plain literal (20), the remaining columns now pull the mast
# height from the "Mast + Tower" sheet (18 m fixed structure + B20/B18).
$headend.Range("F49").Value = 20
$headend.Range("G49").Formula = "=18+'Mast + Tower'!B20"
$headend.Range("H49").Formula = "=18+'Mast + Tower'!B18"
$headend.Range("I49").Formula = "=18+'Mast + Tower'!B18"
$headend.Range("J49").Formula = "=18+'Mast + Tower'!B18"
$headend.Range("K49").Formula = "=18+'Mast + Tower'!B18"
$headend.Range("L49").Formula = "=18+'Mast + Tower'!B18"
$headend.Range("M49").Formula = "=18+'Mast + Tower'!B18"
$headend.Range("N49").Formula = "=18+'Mast + Tower'!B18"

# --- Fix the SUM() range typo in the mast-amplifier S/I check (row 96) -----
$headend.Range("J96").Formula = "=-20*LOG10(10^(-J91/20)+10^(-(J28+2*(J27-7.5*LOG10(SUM(H86:N86)-1)-J32-J63))/20))"

# --- Clear leftover scratch notes that are no longer needed ----------------
$headend.Range("P5").ClearContents()
$headend.Range("F48").ClearContents()
$mast.Range("L16").ClearContents()
